$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move existing values from column B to column A
$ws.Range("A1").Value = $ws.Range("B1").Value2
$ws.Range("A2").Value = $ws.Range("B2").Value2
$ws.Range("A3").Value = $ws.Range("B3").Value2
$ws.Range("A4").Value = $ws.Range("B4").Value2

# Clear old column B data
$ws.Range("B1:B4").Clear()

# Add new row with value 123456
$ws.Range("A5").Value = 123456

# Update selection to match target state
$ws.Range("C8").Select()
